$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Rename sheets "depot" -> "depots" and "customer" -> "customers".
#    Excel automatically rewrites formulas referencing these sheets
#    (e.g. in OLDdistances_depot_costumers) when a sheet is renamed.
# ---------------------------------------------------------------------
$wsDepot = $wb.Worksheets.Item("depot")
$wsDepot.Name = "depots"

$wsCustomer = $wb.Worksheets.Item("customer")
$wsCustomer.Name = "customers"

# ---------------------------------------------------------------------
# 2) Rebuild the "trucks" sheet: insert a new id_truck numeric id
#    column, rename the old id column header to id_type, and append
#    an id_depot column assigning 10 trucks to each of 6 depots.
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("trucks")

# Insert a new column before A; this shifts the old A:E -> B:F and
# keeps all existing data / number formats / column widths intact.
$ws.Columns.Item(1).Insert()

# Fix up the header row.
$ws.Cells.Item(1, 1).Value = "id_truck"
$ws.Cells.Item(1, 2).Value = "id_type"
$ws.Cells.Item(1, 7).Value = "id_depot"

$types = @("Rigido3Ejes", "TrailerTorito", "Trailer3Ejes", "Rigido2Ejes", "Rigido2Ejes")
$weights = @(15000, 24000, 27000, 10000, 10000)
$widths = @(2.5, 2.5, 2.5, 2.5, 2.5)
$lengths = @(9.5, 11, 13.6, 7.5, 7.5)
$heights = @(2.5, 2.5, 2.5, 2.5, 2.5)
$depots = @("F1", "F2", "F3", "F4", "F5", "F6")

$truckId = 1
for ($d = 0; $d -lt 6; $d++) {
    for ($i = 0; $i -lt 10; $i++) {
        $idx = $i % 5
        $row = $truckId + 1

        $ws.Cells.Item($row, 1).Value = $truckId
        $ws.Cells.Item($row, 2).Value = $types[$idx]
        $ws.Cells.Item($row, 3).Value = $weights[$idx]
        $ws.Cells.Item($row, 4).Value = $widths[$idx]
        $ws.Cells.Item($row, 5).Value = $lengths[$idx]
        $ws.Cells.Item($row, 6).Value = $heights[$idx]
        $ws.Cells.Item($row, 7).Value = $depots[$d]

        $truckId = $truckId + 1
    }
}

# ---------------------------------------------------------------------
# 3) Cosmetic selection changes on the depots / customers sheets.
# ---------------------------------------------------------------------
$wsDepot.Range("A2").Select()
$wsCustomer.Range("E29").Select()

# ---------------------------------------------------------------------
# 4) Make "trucks" the active sheet/selection, matching the saved view
#    state in the target workbook.
# ---------------------------------------------------------------------
$ws.Activate()
$ws.Range("G61").Select()
